$wb = $excel.ActiveWorkbook

# Rename the INDUSTRIES sheet to SITES
$ws1 = $wb.Worksheets.Item("ASSESSMENTS")
$ws2 = $wb.Worksheets.Item("INDUSTRIES")
$ws2.Name = "SITES"

# --- ASSESSMENTS sheet: update instruction text (industry -> site wording) ---
$ws1.Range("E5").Value = "Add the assessments you want to add to the tool on the first sheet.`n If the name of any assessment coincides with an existing one, the latter (along with its sites) will be deleted."
$ws1.Range("E6").Value = "In the sites tab, define the industries you want to add,  `nwith the assessment to which it belongs. This assessment does not have to be defined in the assessment sheet, it can be previously defined in the web tool."

# --- SITES sheet (formerly INDUSTRIES): rename headers/labels ---
$ws2.Range("A1").Value = "SITE"
$ws2.Range("E1").Value = "SUB-SUPPLIERS"
$ws2.Range("B3").Value = "Site "
$ws2.Range("E2").Value = "Add as many sub-suppliers (Name, Latitude, Longitude) as needed to the right"

# --- Row height tweak on ASSESSMENTS sheet row 6 ---
$ws1.Rows.Item(6).RowHeight = 58.5

# --- Selections (active cell) to match final saved view state ---
[void]$ws2.Range("B28").Select()
[void]$ws1.Range("E6").Select()
